# Error Calculations and Plots
# Remove two rows that were filtered out of the "missing_data" extract
# ("RM 232" and "SC 92"), which shifts all subsequent rows up by two,
# and move the single missing-value marker in column E from the row
# that is now "SC 5" to the row that is now "SC 101".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (originally row 26).
$ws.Rows.Item(26).Delete()

# After the first delete, "SC 92" (originally row 28) has shifted up to
# row 27 - delete it too.
$ws.Rows.Item(27).Delete()

# "SC 5" (now row 26) gains the imputed value in column E.
$ws.Cells.Item(26, 5).Value = -5

# "SC 101" (now row 27) loses its column E value (becomes the new
# missing-data cell).
$ws.Cells.Item(27, 5).ClearContents()
